$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I18").Value = -0.1227366928773208
$ws.Range("J18").Value = 0.2977747339815092
$ws.Range("K18").Value = -0.0121057452347735
$ws.Range("L18").Value = 2.172310786008958

$ws.Range("I19").Value = 0.2721110030661992
$ws.Range("J19").Value = 0.5975685436877082
$ws.Range("K19").Value = -0.02932312144336056
$ws.Range("L19").Value = 1.816670092425335
